$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 101
$ws.Range("F3").Value = 128
$ws.Range("F4").Value = 644
$ws.Range("F6").Value = 551
$ws.Range("F7").Value = 1541
$ws.Range("F9").Value = 11801
$ws.Range("F12").Value = 124
$ws.Range("F13").Value = 2120
$ws.Range("F14").Value = 901
$ws.Range("F15").Value = 241
$ws.Range("F16").Value = 61
$ws.Range("F17").Value = 236
$ws.Range("F18").Value = 1200
$ws.Range("F19").Value = 175
$ws.Range("F20").Value = 254
$ws.Range("F21").Value = 744
$ws.Range("F22").Value = 663
$ws.Range("F23").Value = 276
$ws.Range("F25").Value = 734
$ws.Range("F26").Value = 3681
$ws.Range("F27").Value = 3681
$ws.Range("F28").Value = 1075
$ws.Range("F29").Value = 819
$ws.Range("F33").Value = 998
$ws.Range("F34").Value = 40
$ws.Range("F35").Value = 74
$ws.Range("F36").Value = 257
$ws.Range("F37").Value = 22
$ws.Range("F39").Value = 18
$ws.Range("F40").Value = 3459
$ws.Range("F41").Value = 4454
$ws.Range("F42").Value = 5492
$ws.Range("F44").Value = 118
$ws.Range("F45").Value = 161
$ws.Range("F46").Value = 276
$ws.Range("F47").Value = 68
$ws.Range("F48").Value = 33
$ws.Range("F49").Value = 4095
$ws.Range("F50").Value = 105

$ws = $wb.Worksheets.Item(2)
$ws.Range("F3").Value = 4159
$ws.Range("F5").Value = 91
$ws.Range("F12").Value = 750
$ws.Range("F21").Value = 12

$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 754
$ws.Range("F3").Value = 424
$ws.Range("F4").Value = 67

$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 754
$ws.Range("F3").Value = 424
$ws.Range("F4").Value = 67
$ws.Range("F5").Value = 101
$ws.Range("F6").Value = 128
$ws.Range("F7").Value = 644
$ws.Range("F9").Value = 551
$ws.Range("F10").Value = 1541
$ws.Range("F11").Value = 11801
$ws.Range("F14").Value = 124
$ws.Range("F15").Value = 2120
$ws.Range("F16").Value = 901
$ws.Range("F17").Value = 61
$ws.Range("F18").Value = 1200
$ws.Range("F19").Value = 175
$ws.Range("F20").Value = 254
$ws.Range("F21").Value = 4159
$ws.Range("F23").Value = 276
$ws.Range("F24").Value = 734
$ws.Range("F25").Value = 3681
$ws.Range("F26").Value = 1075
$ws.Range("F27").Value = 91
$ws.Range("F29").Value = 819
$ws.Range("F31").Value = 998
$ws.Range("F32").Value = 40
$ws.Range("F33").Value = 74
$ws.Range("F34").Value = 257
$ws.Range("F35").Value = 22
$ws.Range("F36").Value = 18
$ws.Range("F37").Value = 4454
$ws.Range("F39").Value = 118
$ws.Range("F40").Value = 161
$ws.Range("F41").Value = 276
$ws.Range("F44").Value = 68
$ws.Range("F45").Value = 33
$ws.Range("F48").Value = 12
$ws.Range("F50").Value = 105
